$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F39").Copy()
$ws.Range("ZZ1").PasteSpecial(-4122)
$ws.Range("F39").Value = "https://youtu.be/bNLYRnJj6R0 "
$ws.Hyperlinks.Add($ws.Range("F39"), "https://youtu.be/bNLYRnJj6R0 ", "", "", "https://youtu.be/bNLYRnJj6R0 ")
$ws.Range("ZZ1").Copy()
$ws.Range("F39").PasteSpecial(-4122)
$ws.Range("ZZ1").Clear()
